# Weekly update: a new week of price records for Membrillo (Mercado Mayorista
# Lo Valledor de Santiago) is inserted at the top of the data block (row 111),
# pushing all the existing data rows down by two rows. The sheet's used range
# therefore grows from A1:T206 to A1:T208.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 111; everything currently at row 111 and
# below (through row 206) shifts down to rows 113-208.
$ws.Rows("111:112").Insert()

# ---- New row 111 ----
$ws.Cells.Item(111, 1).Value2  = 6
$ws.Cells.Item(111, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(111, 3).Value2  = "Metropolitana"
$ws.Cells.Item(111, 4).Value2  = 45079
$ws.Cells.Item(111, 5).Value2  = 13
$ws.Cells.Item(111, 6).Value2  = "Fruta"
$ws.Cells.Item(111, 7).Value2  = 100104
$ws.Cells.Item(111, 8).Value2  = "Frutos de pepita"
$ws.Cells.Item(111, 9).Value2  = 100104003
$ws.Cells.Item(111, 10).Value2 = "Membrillo"
$ws.Cells.Item(111, 11).Value2 = "Champion"
$ws.Cells.Item(111, 12).Value2 = "Especial"
$ws.Cells.Item(111, 13).Value2 = 4
$ws.Cells.Item(111, 14).Value2 = 180000
$ws.Cells.Item(111, 15).Value2 = 180000
$ws.Cells.Item(111, 16).Value2 = 180000
$ws.Cells.Item(111, 17).Value2 = "$/bins (450 kilos)"
$ws.Cells.Item(111, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(111, 19).Value2 = 400
$ws.Cells.Item(111, 20).Value2 = 450

# ---- New row 112 ----
$ws.Cells.Item(112, 1).Value2  = 6
$ws.Cells.Item(112, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(112, 3).Value2  = "Metropolitana"
$ws.Cells.Item(112, 4).Value2  = 45079
$ws.Cells.Item(112, 5).Value2  = 13
$ws.Cells.Item(112, 6).Value2  = "Fruta"
$ws.Cells.Item(112, 7).Value2  = 100104
$ws.Cells.Item(112, 8).Value2  = "Frutos de pepita"
$ws.Cells.Item(112, 9).Value2  = 100104003
$ws.Cells.Item(112, 10).Value2 = "Membrillo"
$ws.Cells.Item(112, 11).Value2 = "Champion"
$ws.Cells.Item(112, 12).Value2 = "Primera"
$ws.Cells.Item(112, 13).Value2 = 6
$ws.Cells.Item(112, 14).Value2 = 150000
$ws.Cells.Item(112, 15).Value2 = 150000
$ws.Cells.Item(112, 16).Value2 = 150000
$ws.Cells.Item(112, 17).Value2 = "$/bins (450 kilos)"
$ws.Cells.Item(112, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(112, 19).Value2 = 333
$ws.Cells.Item(112, 20).Value2 = 450
